$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 82

# Copy the formatting of the previous row's date cell (style "2" - the
# YYYY-MM-DD HH:MM:SS date format) onto the new date cell, then set values.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 45884
$ws.Cells.Item($row, 2).Value = 0.06654624964350926
